# Card2: "إضافة حدث جديد" — append a new service-event row.
#
# The sheet stores every cell as literal text (even things that look like
# numbers, e.g. "2", "55", "1111"), and blank cells are stored as empty
# text rather than truly-empty cells. To reproduce that on this engine we
# assign values with a leading apostrophe (forces text, Excel strips the
# apostrophe itself) and then reset the cell style to "Normal" so the
# auto-added quote-prefix formatting doesn't stick around.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

function Set-TextCell {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# 1) The previously-blank cells of row 16 become the literal text "nan"
#    (matching every other already-serviced row in this sheet).
Set-TextCell 16 2  "nan"   # B16
Set-TextCell 16 3  "nan"   # C16
Set-TextCell 16 5  "nan"   # E16
Set-TextCell 16 6  "nan"   # F16
Set-TextCell 16 7  "nan"   # G16
Set-TextCell 16 8  "nan"   # H16
Set-TextCell 16 9  "nan"   # I16
Set-TextCell 16 10 "nan"   # J16
Set-TextCell 16 11 "nan"   # K16
Set-TextCell 16 13 "nan"   # M16
Set-TextCell 16 14 "nan"   # N16

# 2) A brand-new row 17 is appended, carrying the same "new event" template
#    row 16 originally had (before it got filled in above): card "2",
#    Tones "55", Date "1111", Serviced by "ححح", everything else blank.
Set-TextCell 17 1  "2"     # A17
Set-TextCell 17 2  ""      # B17
Set-TextCell 17 3  ""      # C17
Set-TextCell 17 4  "55"    # D17
Set-TextCell 17 5  ""      # E17
Set-TextCell 17 6  ""      # F17
Set-TextCell 17 7  ""      # G17
Set-TextCell 17 8  ""      # H17
Set-TextCell 17 9  ""      # I17
Set-TextCell 17 10 ""      # J17
Set-TextCell 17 11 ""      # K17
Set-TextCell 17 12 "1111"  # L17
Set-TextCell 17 13 ""      # M17
Set-TextCell 17 14 ""      # N17
Set-TextCell 17 15 "ححح"  # O17
